$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved for all edited cells (data is inline text,
# not numeric, even though some values look like numbers e.g. "596.74").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.078.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.40"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.68%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.120.33"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.475"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.630.30"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.284.36"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.113"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.136.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.71"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.84"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.04"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.66"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.60"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0739"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.26"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -10.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0400"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.877.22"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.273"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.77"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.11%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.115"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.13"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.53%  "
